# Scen_NCAP_NUC.xlsx - "Add files via upload" commit
# Replaces the old single MIN_PINK_HYDROGEN lower-bound row (row 45) with a
# new set of LO / ACT_BND rows for PRE_HYDROGEN spanning 2030-2050.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45: was UP/ACT_BND/2035/250/MIN_PINK_HYDROGEN -------------------
# becomes LO/ACT_BND/2030/20/PRE_HYDROGEN
$ws.Range("B45").Value = "LO"
$ws.Range("C45").Value = "ACT_BND"
$ws.Range("D45").Value = 2030
$ws.Range("E45").Value = 20
$ws.Range("F45").ClearFormats()
$ws.Range("F45").Value = "PRE_HYDROGEN"

# --- New rows 46-49: LO / ACT_BND / PRE_HYDROGEN bound series -------------
$ws.Range("B46").Value = "LO"
$ws.Range("C46").Value = "ACT_BND"
$ws.Range("D46").Value = 2035
$ws.Range("E46").Value = 80
$ws.Range("F46").Value = "PRE_HYDROGEN"

$ws.Range("B47").Value = "LO"
$ws.Range("C47").Value = "ACT_BND"
$ws.Range("D47").Value = 2040
$ws.Range("E47").Value = 150
$ws.Range("F47").Value = "PRE_HYDROGEN"

$ws.Range("B48").Value = "LO"
$ws.Range("C48").Value = "ACT_BND"
$ws.Range("D48").Value = 2045
$ws.Range("E48").Value = 250
$ws.Range("F48").Value = "PRE_HYDROGEN"

$ws.Range("B49").Value = "LO"
$ws.Range("C49").Value = "ACT_BND"
$ws.Range("D49").Value = 2050
$ws.Range("E49").Value = 250
$ws.Range("F49").Value = "PRE_HYDROGEN"

# --- Sheet view bookkeeping (matches the author's last on-screen state) ---
$ws.Range("J43").Select() | Out-Null
